$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.448.06"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.868.86"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2801"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06509"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "1.840.31"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07446"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.085"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6421"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").Value = "30.429.22"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007492"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "230.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.57%  "
$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.150"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.096"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "170.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.343"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.905"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.20%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1044"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.93%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.385"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.25%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.275"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.985"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04983"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.179"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7406"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("B35").Value = "Frax"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9992"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.710"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01933"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.78%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.631"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9175"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.048"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "105.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9964"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4193"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.579"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.61%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.215"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1227"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.84%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.921"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.423"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05651"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.21%  "
